# ------------------------------------------------------------------
# Applies the "1FC journey" sanity-sheet edit described in the commit:
#   - bump the UserID test-data value on OneFCFlowTestData (B2: data1 -> data5)
#   - add a new "OneFCFlowTestDataSanity" sheet after OneFCFlowTestData
#     with the 1FC (login/dashboard/choose-agreement) test data rows
#   - make the new sheet the active/selected tab
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. OneFCFlowTestData: update the UserID data value used by row 2 ---
$flowSheet = $wb.Worksheets.Item("OneFCFlowTestData")
$flowSheet.Range("B2").Value = "data5"

# --- 2. Add the new sheet right after OneFCFlowTestData ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sanitySheet = $wb.Worksheets.Add($null, $lastSheet)
$sanitySheet.Name = "OneFCFlowTestDataSanity"

# --- 3. Populate the 1FC sanity test data grid (header + 3 data rows) ---
$sanitySheet.Range("A1").Value = "TDID1"
$sanitySheet.Range("B1").Value = "UserID1"
$sanitySheet.Range("C1").Value = "LoginPage"
$sanitySheet.Range("D1").Value = "Dashboard"
$sanitySheet.Range("E1").Value = "Choose a commercial agreement"
$sanitySheet.Range("F1").Value = "Choose a commercial agreement"
$sanitySheet.Range("G1").Value = "Procurement overview"
$sanitySheet.Range("H1").Value = "Procurement overview"
$sanitySheet.Range("I1").Value = "Choose how to find a supplier"
$sanitySheet.Range("J1").Value = "Choose how to find a supplier"
$sanitySheet.Range("K1").Value = "Write and publish your requirements"

$sanitySheet.Range("A2").Value = "TDID"
$sanitySheet.Range("B2").Value = "UserID"
$sanitySheet.Range("C2").Value = "Login_Title"
$sanitySheet.Range("D2").Value = "Dashboard_Title"
$sanitySheet.Range("E2").Value = "ChooseAgreement_Title"
$sanitySheet.Range("F2").Value = "ChooseAgreement_Lot"
$sanitySheet.Range("G2").Value = "S1_Title"
$sanitySheet.Range("H2").Value = "S1_Section"
$sanitySheet.Range("I2").Value = "S2_Title"
$sanitySheet.Range("J2").Value = "S2_Radiobutton"
$sanitySheet.Range("K2").Value = "S3_Title"

$sanitySheet.Range("A3").Value = "TD001"
$sanitySheet.Range("B3").Value = "data1"
$sanitySheet.Range("C3").Value = "Sign in to the Public Procurement Gateway"
$sanitySheet.Range("D3").Value = "Find suppliers and run your procurement online."
$sanitySheet.Range("E3").Value = "Choose a commercial agreement"
$sanitySheet.Range("F3").Value = "Lot 1: Digital Programmes"
$sanitySheet.Range("G3").Value = "Procurement overview"
$sanitySheet.Range("H3").Value = "3. Write and publish your requirements"
$sanitySheet.Range("I3").Value = "Choose how to find a supplier"
$sanitySheet.Range("J3").Value = "1FC"
$sanitySheet.Range("K3").Value = "Write and publish your requirements"

$sanitySheet.Range("A4").Value = "TD002"
$sanitySheet.Range("B4").Value = "data2"
$sanitySheet.Range("C4").Value = "Sign in to the Public Procurement Gateway"
$sanitySheet.Range("D4").Value = "Find suppliers and run your procurement online."
$sanitySheet.Range("E4").Value = "Choose a commercial agreement"
$sanitySheet.Range("F4").Value = "Lot 1: Digital Programmes"
$sanitySheet.Range("G4").Value = "Procurement overview"
$sanitySheet.Range("H4").Value = "3. Write and publish your requirements"
$sanitySheet.Range("I4").Value = "Choose how to find a supplier"
$sanitySheet.Range("J4").Value = "1FC"
$sanitySheet.Range("K4").Value = "Write and publish your requirements"

# --- 4. Formatting: vertical-top alignment on the used grid, wrap text on K3:K4 ---
$sanitySheet.Range("A1:K4").VerticalAlignment = -4160
$sanitySheet.Range("K3:K4").WrapText = $true

# --- 5. Column widths: autofit the text columns, keep the narrow index columns ---
$sanitySheet.Columns.Item(1).ColumnWidth = 8.7265625
$sanitySheet.Columns.Item(2).ColumnWidth = 8.7265625
$sanitySheet.Range("A1:K4").Columns.AutoFit()
$sanitySheet.Columns.Item(12).ColumnWidth = 18.7265625
$sanitySheet.Columns.Item(13).ColumnWidth = 30.453125
$sanitySheet.Columns.Item(14).ColumnWidth = 30.453125
$sanitySheet.Columns.Item(15).ColumnWidth = 44.1796875

# --- 6. Page setup to match the sibling sheets ---
$sanitySheet.PageSetup.Orientation = 1

# --- 7. Selections: leave OneFCFlowTestData's cursor at A2 (no longer the active tab) ---
$flowSheet.Activate()
$flowSheet.Range("A2").Select()

# --- 8. Activate the new sheet and select its working cell, making it the active tab ---
$sanitySheet.Activate()
$sanitySheet.Range("I6").Select()
